$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 990.5714
$ws.Range("I33").Value = 999.94446
$ws.Range("J33").Value = 934.3333
$ws.Range("K33").Value = 999.94446
$ws.Range("L33").Value = 934.3333
$ws.Range("M33").Value = -770.94446
$ws.Range("N33").Value = -1392.3333
$ws.Range("H41").Value = 26319562
$ws.Range("I41").Value = 787.5
$ws.Range("K41").Value = 787.5
$ws.Range("M41").Value = -347.5
$ws.Range("H62").Value = 25004998
$ws.Range("I62").Value = 31253748
$ws.Range("J62").Value = 9999
$ws.Range("K62").Value = 31253748
$ws.Range("L62").Value = 9999
$ws.Range("M62").Value = -31253124
$ws.Range("N62").Value = -11247
$ws.Range("H65").Value = 25004998
$ws.Range("I65").Value = 31253748
$ws.Range("J65").Value = 9999
$ws.Range("K65").Value = 156268740
$ws.Range("L65").Value = 49995
$ws.Range("M65").Value = -156265620
$ws.Range("N65").Value = -56235
$ws.Range("H96").Value = 473.53333
$ws.Range("I96").Value = 385.6
$ws.Range("K96").Value = 1156.8
$ws.Range("M96").Value = 216.1999999999998
$ws.Range("H106").Value = 3249.875
$ws.Range("H107").Value = 35719.516
$ws.Range("I107").Value = 44329.26
$ws.Range("K107").Value = 44329.26
$ws.Range("M107").Value = -42409.26
$ws.Range("H112").Value = 4100.4136
$ws.Range("J112").Value = 4524.48
$ws.Range("L112").Value = 13573.44
$ws.Range("N112").Value = -15789.44
$ws.Range("H116").Value = 6407.091
$ws.Range("I116").Value = 6068.2856
$ws.Range("K116").Value = 6068.2856
$ws.Range("M116").Value = -2626.2856
$ws.Range("H137").Value = 5303.9375
$ws.Range("I137").Value = 3775.7
$ws.Range("K137").Value = 11327.1
$ws.Range("M137").Value = -8777.099999999999
$ws.Range("H138").Value = 6139.3057
$ws.Range("I138").Value = 4761.6875
$ws.Range("J138").Value = 6532.9106
$ws.Range("K138").Value = 14285.0625
$ws.Range("L138").Value = 19598.7318
$ws.Range("M138").Value = -9145.0625
$ws.Range("N138").Value = -29878.7318

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3589.5476
$ws.Range("I32").Value = 2737.5386
$ws.Range("K32").Value = 2737.5386
$ws.Range("M32").Value = -2450.5386
$ws.Range("H45").Value = 3906.68
$ws.Range("I45").Value = 2832.4375
$ws.Range("J45").Value = 5816.4443
$ws.Range("K45").Value = 2832.4375
$ws.Range("L45").Value = 5816.4443
$ws.Range("M45").Value = -2455.4375
$ws.Range("N45").Value = -6570.4443
$ws.Range("H132").Value = 5879.6665
$ws.Range("I132").Value = 4098.3335
$ws.Range("K132").Value = 12295.0005
$ws.Range("M132").Value = -9765.000499999998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1216800.9
$ws.Range("I86").Value = 1891032.9
$ws.Range("J86").Value = 3183.4
$ws.Range("K86").Value = 1891032.9
$ws.Range("L86").Value = 3183.4
$ws.Range("M86").Value = -1889909.9
$ws.Range("N86").Value = -5429.4
$ws.Range("H89").Value = 1216800.9
$ws.Range("I89").Value = 1891032.9
$ws.Range("J89").Value = 3183.4
$ws.Range("K89").Value = 9455164.5
$ws.Range("L89").Value = 15917
$ws.Range("M89").Value = -9449548.5
$ws.Range("N89").Value = -27149
$ws.Range("H134").Value = 14575.702
$ws.Range("I134").Value = 1706.1846
$ws.Range("K134").Value = 5118.5538
$ws.Range("M134").Value = -2583.5538
$ws.Range("H140").Value = 100000
$ws.Range("J140").Value = 100000
$ws.Range("L140").Value = 100000
$ws.Range("N140").Value = -110360

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 545.26666
$ws.Range("I22").Value = 318.5
$ws.Range("K22").Value = 318.5
$ws.Range("M22").Value = 31.5
$ws.Range("H31").Value = 720905.6
$ws.Range("I31").Value = 1001968
$ws.Range("K31").Value = 1001968
$ws.Range("M31").Value = -1001673
$ws.Range("H34").Value = 720905.6
$ws.Range("I34").Value = 1001968
$ws.Range("K34").Value = 1001968
$ws.Range("M34").Value = -1001766
$ws.Range("H139").Value = 97942.5
$ws.Range("J139").Value = 98385
$ws.Range("L139").Value = 98385
$ws.Range("N139").Value = -108665

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value = 1585
$ws.Range("I139").Value = 881.6667
$ws.Range("J139").Value = 4750
$ws.Range("K139").Value = 2645.0001
$ws.Range("L139").Value = 14250
$ws.Range("M139").Value = 2494.9999
$ws.Range("N139").Value = -24530

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H22").Value = 8649.799999999999
$ws.Range("I22").Value = 1620
$ws.Range("J22").Value = 13336.333
$ws.Range("K22").Value = 1620
$ws.Range("L22").Value = 13336.333
$ws.Range("M22").Value = -1091
$ws.Range("N22").Value = -14394.333
$ws.Range("H113").Value = 386477.28
$ws.Range("I113").Value = 589799.9399999999
$ws.Range("J113").Value = 2423.3333
$ws.Range("K113").Value = 589799.9399999999
$ws.Range("L113").Value = 2423.3333
$ws.Range("M113").Value = -587629.9399999999
$ws.Range("N113").Value = -6763.3333
$ws.Range("H132").Value = 157362.55
$ws.Range("I132").Value = 187340.67
$ws.Range("K132").Value = 562022.01
$ws.Range("M132").Value = -559492.01

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2060.1
$ws.Range("I82").Value = 2434.8333
$ws.Range("K82").Value = 2434.8333
$ws.Range("M82").Value = -2073.8333
$ws.Range("H85").Value = 2060.1
$ws.Range("I85").Value = 2434.8333
$ws.Range("K85").Value = 2434.8333
$ws.Range("M85").Value = -1186.8333
$ws.Range("H88").Value = 15000
$ws.Range("J88").Value = 15000
$ws.Range("L88").Value = 15000
$ws.Range("N88").Value = -15856
$ws.Range("H91").Value = 15000
$ws.Range("J91").Value = 15000
$ws.Range("L91").Value = 15000
$ws.Range("N91").Value = -17964
$ws.Range("H136").Value = 2961.2285
$ws.Range("I136").Value = 2257
$ws.Range("J136").Value = 4311
$ws.Range("K136").Value = 6771
$ws.Range("L136").Value = 12933
$ws.Range("M136").Value = -4221
$ws.Range("N136").Value = -18033

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H44").Value = 6644.4443
$ws.Range("J44").Value = 6644.4443
$ws.Range("L44").Value = 6644.4443
$ws.Range("N44").Value = -7752.4443
$ws.Range("H120").Value = 82494.5
$ws.Range("J120").Value = 82494.5
$ws.Range("L120").Value = 82494.5
$ws.Range("N120").Value = -92170.5
